$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:r='http://schemas.openxmlformats.org/officeDocument/2006/relationships'"

# ---------------------------------------------------------------------------
# Hunk 1: merge the "Bình Định, ngày ..." paragraph with the following empty
# paragraph that carries the section break, enlarging the date-line font,
# right-aligning it, and giving the section a two explicit-width column
# layout.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p6 = $d.Paragraphs.Item(6)
$mergedRange = $d.Range($p5.Range.Start, $p6.Range.End)
$mergedXml = @"
<w:p $wNs>
  <w:pPr>
    <w:pBdr>
      <w:top w:val="nil"/>
      <w:left w:val="nil"/>
      <w:bottom w:val="nil"/>
      <w:right w:val="nil"/>
      <w:between w:val="nil"/>
    </w:pBdr>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:jc w:val="right"/>
    <w:rPr>
      <w:color w:val="000000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:sectPr>
      <w:footerReference w:type="default" r:id="rId7"/>
      <w:pgSz w:w="11907" w:h="16840"/>
      <w:pgMar w:top="1134" w:right="1021" w:bottom="1077" w:left="1701" w:header="567" w:footer="567" w:gutter="0"/>
      <w:pgNumType w:start="1"/>
      <w:cols w:num="2" w:space="170" w:equalWidth="0">
        <w:col w:w="4082" w:space="170"/>
        <w:col w:w="4933"/>
      </w:cols>
    </w:sectPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve">Bình Định, </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:color w:val="000000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>ngày {Day} tháng {Month} năm {Year}</w:t>
  </w:r>
</w:p>
"@
$mergedRange.InsertXML($mergedXml)

Write-Output "hunk1 done"

# ---------------------------------------------------------------------------
# Hunk 2: "Điều 3. ..." heading paragraph - spacing before/after 60 -> 120,
# drop the firstLine indent.
# ---------------------------------------------------------------------------
$pDieu3 = $d.Paragraphs.Item(27)
$dieu3Xml = @"
<w:p $wNs>
  <w:pPr>
    <w:pBdr>
      <w:top w:val="nil"/>
      <w:left w:val="nil"/>
      <w:bottom w:val="nil"/>
      <w:right w:val="nil"/>
      <w:between w:val="nil"/>
    </w:pBdr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="120" w:after="120"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:b/>
      <w:color w:val="000000"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t>Điều 3. Các quyền lợi, quyền hạn và nghĩa vụ của Bên B:</w:t>
  </w:r>
</w:p>
"@
$pDieu3.Range.InsertXML($dieu3Xml)

Write-Output "hunk2 done"

# ---------------------------------------------------------------------------
# Hunk 3 & 4: move <w:lastRenderedPageBreak/> from the start of the
# "- Bên B có quyền đề xuất..." run to the end of the "3.2 Quyền hạn: " run.
# ---------------------------------------------------------------------------
$pQuyenHan = $d.Paragraphs.Item(32)
$quyenHanXml = @"
<w:p $wNs>
  <w:pPr>
    <w:pBdr>
      <w:top w:val="nil"/>
      <w:left w:val="nil"/>
      <w:bottom w:val="nil"/>
      <w:right w:val="nil"/>
      <w:between w:val="nil"/>
    </w:pBdr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="3969"/>
    </w:tabs>
    <w:spacing w:before="120" w:after="120"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:b/>
      <w:color w:val="000000"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">3.2 Quyền hạn: </w:t>
  </w:r>
</w:p>
"@
$pQuyenHan.Range.InsertXML($quyenHanXml)

$pBenB = $d.Paragraphs.Item(33)
$benBXml = @"
<w:p $wNs>
  <w:pPr>
    <w:pBdr>
      <w:top w:val="nil"/>
      <w:left w:val="nil"/>
      <w:bottom w:val="nil"/>
      <w:right w:val="nil"/>
      <w:between w:val="nil"/>
    </w:pBdr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:tabs>
      <w:tab w:val="left" w:pos="8931"/>
    </w:tabs>
    <w:spacing w:before="60" w:after="60"/>
    <w:ind w:firstLine="709"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:color w:val="000000"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t xml:space="preserve">- Bên B có quyền đề xuất, khiếu nại, thay đổi, tạm hoãn, chấm dứt hợp đồng theo quy định của </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t>p</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t>háp luật;</w:t>
  </w:r>
</w:p>
"@
$pBenB.Range.InsertXML($benBXml)

Write-Output "hunk3+4 done"

# ---------------------------------------------------------------------------
# Hunk 5 & 6: split the confidentiality-obligations paragraph into two runs,
# moving <w:lastRenderedPageBreak/> to just before "giải mã, ..." and
# removing it from the following paragraph ("- Thông báo ngay lập tức...").
# ---------------------------------------------------------------------------
$pConfid = $d.Paragraphs.Item(48)
$confidXml = @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="120" w:line="320" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">- Có trách nhiệm không tiết lộ Thông tin bảo mật nào cho bên thứ ba, ngoại trừ trường hợp được ủy quyền bằng văn bản của bên A; không sử dụng hoặc khai thác Thông tin bảo mật của bên A, ngoại trừ phục vụ cho mục đích được cho phép; không sao chép, mô phỏng, không dùng kỹ thuật phân tích ngược, không phân tích, lắp đặt, </w:t>
  </w:r>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>giải mã, phân phối hoặc chuyển giao bất kỳ Thông tin bảo mật nào của bên A, trừ khi được ủy quyền bằng văn bản của bên A.</w:t>
  </w:r>
</w:p>
"@
$pConfid.Range.InsertXML($confidXml)

$pThongBao = $d.Paragraphs.Item(49)
$thongBaoXml = @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="120" w:line="320" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>- Thông báo ngay lập tức cho bên A trong trường hợp có bất kỳ sự thất thoát hoặc tiết lộ trái phép Thông tin bảo mật của bên A.</w:t>
  </w:r>
</w:p>
"@
$pThongBao.Range.InsertXML($thongBaoXml)

Write-Output "hunk5+6 done"
